# ===========================================================================
# 688559-海目星.xlsx -- add "2022-Q3" quarterly fund-holders sheet
# ===========================================================================
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q3" worksheet right before "2022-Q2"
# ---------------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($refSheet)
$q3.Name = "2022-Q3"

# Header row (bold, centered, thin border - matches the other quarter sheets)
$q3Header = $q3.Range("B1:H1")
$q3Header.Font.Bold = $true
$q3Header.HorizontalAlignment = -4108
$q3Header.VerticalAlignment = -4160
$q3Header.Borders.LineStyle = 1

$q3.Cells.Item(1, 2).Value = "基金代码"
$q3.Cells.Item(1, 3).Value = "基金名称"
$q3.Cells.Item(1, 4).Value = "基金规模"
$q3.Cells.Item(1, 5).Value = "股票总仓位"
$q3.Cells.Item(1, 6).Value = "仓位占比"
$q3.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q3.Cells.Item(1, 8).Value = "仓位排名"

# Columns B:G hold text-like data (fund codes with leading zeros, and
# numbers formatted as fixed-decimal strings) - force text storage so
# values such as "010122" or "21.73" are not coerced into numbers.
$q3.Range("B2:G45").NumberFormat = "@"

# ID column A (bold, centered, thin border) + data rows 2-45
$q3ColA = $q3.Range("A2:A45")
$q3ColA.Font.Bold = $true
$q3ColA.HorizontalAlignment = -4108
$q3ColA.VerticalAlignment = -4160
$q3ColA.Borders.LineStyle = 1

$q3.Cells.Item(2, 1).Value = 0
$q3.Cells.Item(2, 2).Value = "398021"
$q3.Cells.Item(2, 3).Value = "中海能源策略混合"
$q3.Cells.Item(2, 4).Value = "21.73"
$q3.Cells.Item(2, 5).Value = "88.03"
$q3.Cells.Item(2, 6).Value = "4.70"
$q3.Cells.Item(2, 7).Value = "1.0213"
$q3.Cells.Item(2, 8).Value = 4

$q3.Cells.Item(3, 1).Value = 1
$q3.Cells.Item(3, 2).Value = "202027"
$q3.Cells.Item(3, 3).Value = "南方高端装备灵活配置混合A"
$q3.Cells.Item(3, 4).Value = "12.54"
$q3.Cells.Item(3, 5).Value = "93.38"
$q3.Cells.Item(3, 6).Value = "4.61"
$q3.Cells.Item(3, 7).Value = "0.5781"
$q3.Cells.Item(3, 8).Value = 4

$q3.Cells.Item(4, 1).Value = 2
$q3.Cells.Item(4, 2).Value = "460001"
$q3.Cells.Item(4, 3).Value = "华泰柏瑞盛世中国混合"
$q3.Cells.Item(4, 4).Value = "18.02"
$q3.Cells.Item(4, 5).Value = "89.00"
$q3.Cells.Item(4, 6).Value = "2.89"
$q3.Cells.Item(4, 7).Value = "0.5208"
$q3.Cells.Item(4, 8).Value = 10

$q3.Cells.Item(5, 1).Value = 3
$q3.Cells.Item(5, 2).Value = "501186"
$q3.Cells.Item(5, 3).Value = "华夏兴融灵活配置混合（LOF）A"
$q3.Cells.Item(5, 4).Value = "8.78"
$q3.Cells.Item(5, 5).Value = "89.61"
$q3.Cells.Item(5, 6).Value = "2.24"
$q3.Cells.Item(5, 7).Value = "0.1967"
$q3.Cells.Item(5, 8).Value = 10

$q3.Cells.Item(6, 1).Value = 4
$q3.Cells.Item(6, 2).Value = "010122"
$q3.Cells.Item(6, 3).Value = "华泰柏瑞优势领航混合A"
$q3.Cells.Item(6, 4).Value = "3.95"
$q3.Cells.Item(6, 5).Value = "79.94"
$q3.Cells.Item(6, 6).Value = "3.52"
$q3.Cells.Item(6, 7).Value = "0.1390"
$q3.Cells.Item(6, 8).Value = 5

$q3.Cells.Item(7, 1).Value = 5
$q3.Cells.Item(7, 2).Value = "000082"
$q3.Cells.Item(7, 3).Value = "嘉实研究阿尔法股票"
$q3.Cells.Item(7, 4).Value = "6.18"
$q3.Cells.Item(7, 5).Value = "90.88"
$q3.Cells.Item(7, 6).Value = "1.88"
$q3.Cells.Item(7, 7).Value = "0.1162"
$q3.Cells.Item(7, 8).Value = 8

$q3.Cells.Item(8, 1).Value = 6
$q3.Cells.Item(8, 2).Value = "004895"
$q3.Cells.Item(8, 3).Value = "华商鑫安灵活配置混合"
$q3.Cells.Item(8, 4).Value = "2.11"
$q3.Cells.Item(8, 5).Value = "92.54"
$q3.Cells.Item(8, 6).Value = "4.21"
$q3.Cells.Item(8, 7).Value = "0.0888"
$q3.Cells.Item(8, 8).Value = 6

$q3.Cells.Item(9, 1).Value = 7
$q3.Cells.Item(9, 2).Value = "005207"
$q3.Cells.Item(9, 3).Value = "南方高端装备灵活配置混合C"
$q3.Cells.Item(9, 4).Value = "1.83"
$q3.Cells.Item(9, 5).Value = "93.38"
$q3.Cells.Item(9, 6).Value = "4.61"
$q3.Cells.Item(9, 7).Value = "0.0844"
$q3.Cells.Item(9, 8).Value = 4

$q3.Cells.Item(10, 1).Value = 8
$q3.Cells.Item(10, 2).Value = "960019"
$q3.Cells.Item(10, 3).Value = "招商行業領先混合型證券投資基金 H"
$q3.Cells.Item(10, 4).Value = "2.80"
$q3.Cells.Item(10, 5).Value = "89.41"
$q3.Cells.Item(10, 6).Value = "2.96"
$q3.Cells.Item(10, 7).Value = "0.0829"
$q3.Cells.Item(10, 8).Value = 10

$q3.Cells.Item(11, 1).Value = 9
$q3.Cells.Item(11, 2).Value = "217012"
$q3.Cells.Item(11, 3).Value = "招商行业领先混合A"
$q3.Cells.Item(11, 4).Value = "2.80"
$q3.Cells.Item(11, 5).Value = "89.41"
$q3.Cells.Item(11, 6).Value = "2.96"
$q3.Cells.Item(11, 7).Value = "0.0829"
$q3.Cells.Item(11, 8).Value = 10

$q3.Cells.Item(12, 1).Value = 10
$q3.Cells.Item(12, 2).Value = "580006"
$q3.Cells.Item(12, 3).Value = "东吴新经济混合A"
$q3.Cells.Item(12, 4).Value = "1.37"
$q3.Cells.Item(12, 5).Value = "90.50"
$q3.Cells.Item(12, 6).Value = "5.66"
$q3.Cells.Item(12, 7).Value = "0.0775"
$q3.Cells.Item(12, 8).Value = 7

$q3.Cells.Item(13, 1).Value = 11
$q3.Cells.Item(13, 2).Value = "001815"
$q3.Cells.Item(13, 3).Value = "华泰柏瑞激励动力灵活配置混合A"
$q3.Cells.Item(13, 4).Value = "1.93"
$q3.Cells.Item(13, 5).Value = "67.38"
$q3.Cells.Item(13, 6).Value = "3.51"
$q3.Cells.Item(13, 7).Value = "0.0677"
$q3.Cells.Item(13, 8).Value = 2

$q3.Cells.Item(14, 1).Value = 12
$q3.Cells.Item(14, 2).Value = "001808"
$q3.Cells.Item(14, 3).Value = "银华互联网主题灵活配置混合A"
$q3.Cells.Item(14, 4).Value = "1.83"
$q3.Cells.Item(14, 5).Value = "89.52"
$q3.Cells.Item(14, 6).Value = "3.61"
$q3.Cells.Item(14, 7).Value = "0.0661"
$q3.Cells.Item(14, 8).Value = 10

$q3.Cells.Item(15, 1).Value = 13
$q3.Cells.Item(15, 2).Value = "010238"
$q3.Cells.Item(15, 3).Value = "安信创新先锋混合C"
$q3.Cells.Item(15, 4).Value = "1.10"
$q3.Cells.Item(15, 5).Value = "93.80"
$q3.Cells.Item(15, 6).Value = "5.90"
$q3.Cells.Item(15, 7).Value = "0.0649"
$q3.Cells.Item(15, 8).Value = 3

$q3.Cells.Item(16, 1).Value = 14
$q3.Cells.Item(16, 2).Value = "002707"
$q3.Cells.Item(16, 3).Value = "摩根士丹利华鑫科技领先灵活配置混合A"
$q3.Cells.Item(16, 4).Value = "1.75"
$q3.Cells.Item(16, 5).Value = "94.13"
$q3.Cells.Item(16, 6).Value = "3.18"
$q3.Cells.Item(16, 7).Value = "0.0556"
$q3.Cells.Item(16, 8).Value = 9

$q3.Cells.Item(17, 1).Value = 15
$q3.Cells.Item(17, 2).Value = "004044"
$q3.Cells.Item(17, 3).Value = "金鹰转型动力灵活配置混合"
$q3.Cells.Item(17, 4).Value = "0.60"
$q3.Cells.Item(17, 5).Value = "88.13"
$q3.Cells.Item(17, 6).Value = "7.34"
$q3.Cells.Item(17, 7).Value = "0.0440"
$q3.Cells.Item(17, 8).Value = 3

$q3.Cells.Item(18, 1).Value = 16
$q3.Cells.Item(18, 2).Value = "002082"
$q3.Cells.Item(18, 3).Value = "华泰柏瑞激励动力灵活配置混合C"
$q3.Cells.Item(18, 4).Value = "0.98"
$q3.Cells.Item(18, 5).Value = "67.38"
$q3.Cells.Item(18, 6).Value = "3.51"
$q3.Cells.Item(18, 7).Value = "0.0344"
$q3.Cells.Item(18, 8).Value = 2

$q3.Cells.Item(19, 1).Value = 17
$q3.Cells.Item(19, 2).Value = "673090"
$q3.Cells.Item(19, 3).Value = "西部利得个股精选股票A"
$q3.Cells.Item(19, 4).Value = "1.25"
$q3.Cells.Item(19, 5).Value = "86.69"
$q3.Cells.Item(19, 6).Value = "2.65"
$q3.Cells.Item(19, 7).Value = "0.0331"
$q3.Cells.Item(19, 8).Value = 3

$q3.Cells.Item(20, 1).Value = 18
$q3.Cells.Item(20, 2).Value = "012617"
$q3.Cells.Item(20, 3).Value = "东吴新经济混合C"
$q3.Cells.Item(20, 4).Value = "0.55"
$q3.Cells.Item(20, 5).Value = "90.50"
$q3.Cells.Item(20, 6).Value = "5.66"
$q3.Cells.Item(20, 7).Value = "0.0311"
$q3.Cells.Item(20, 8).Value = 7

$q3.Cells.Item(21, 1).Value = 19
$q3.Cells.Item(21, 2).Value = "006547"
$q3.Cells.Item(21, 3).Value = "红塔红土盛弘灵活配置混合A"
$q3.Cells.Item(21, 4).Value = "0.98"
$q3.Cells.Item(21, 5).Value = "86.31"
$q3.Cells.Item(21, 6).Value = "3.15"
$q3.Cells.Item(21, 7).Value = "0.0309"
$q3.Cells.Item(21, 8).Value = 8

$q3.Cells.Item(22, 1).Value = 20
$q3.Cells.Item(22, 2).Value = "013262"
$q3.Cells.Item(22, 3).Value = "西部利得个股精选股票C"
$q3.Cells.Item(22, 4).Value = "0.95"
$q3.Cells.Item(22, 5).Value = "86.69"
$q3.Cells.Item(22, 6).Value = "2.65"
$q3.Cells.Item(22, 7).Value = "0.0252"
$q3.Cells.Item(22, 8).Value = 3

$q3.Cells.Item(23, 1).Value = 21
$q3.Cells.Item(23, 2).Value = "010237"
$q3.Cells.Item(23, 3).Value = "安信创新先锋混合A"
$q3.Cells.Item(23, 4).Value = "0.41"
$q3.Cells.Item(23, 5).Value = "93.80"
$q3.Cells.Item(23, 6).Value = "5.90"
$q3.Cells.Item(23, 7).Value = "0.0242"
$q3.Cells.Item(23, 8).Value = 3

$q3.Cells.Item(24, 1).Value = 22
$q3.Cells.Item(24, 2).Value = "014360"
$q3.Cells.Item(24, 3).Value = "红塔红土稳健添利混合A"
$q3.Cells.Item(24, 4).Value = "1.23"
$q3.Cells.Item(24, 5).Value = "37.80"
$q3.Cells.Item(24, 6).Value = "1.45"
$q3.Cells.Item(24, 7).Value = "0.0178"
$q3.Cells.Item(24, 8).Value = 7

$q3.Cells.Item(25, 1).Value = 23
$q3.Cells.Item(25, 2).Value = "010123"
$q3.Cells.Item(25, 3).Value = "华泰柏瑞优势领航混合C"
$q3.Cells.Item(25, 4).Value = "0.39"
$q3.Cells.Item(25, 5).Value = "79.94"
$q3.Cells.Item(25, 6).Value = "3.52"
$q3.Cells.Item(25, 7).Value = "0.0137"
$q3.Cells.Item(25, 8).Value = 5

$q3.Cells.Item(26, 1).Value = 24
$q3.Cells.Item(26, 2).Value = "006548"
$q3.Cells.Item(26, 3).Value = "红塔红土盛弘灵活配置混合C"
$q3.Cells.Item(26, 4).Value = "0.41"
$q3.Cells.Item(26, 5).Value = "86.31"
$q3.Cells.Item(26, 6).Value = "3.15"
$q3.Cells.Item(26, 7).Value = "0.0129"
$q3.Cells.Item(26, 8).Value = 8

$q3.Cells.Item(27, 1).Value = 25
$q3.Cells.Item(27, 2).Value = "005437"
$q3.Cells.Item(27, 3).Value = "易方达易百智能量化策略灵活配置混合A"
$q3.Cells.Item(27, 4).Value = "0.81"
$q3.Cells.Item(27, 5).Value = "94.77"
$q3.Cells.Item(27, 6).Value = "1.09"
$q3.Cells.Item(27, 7).Value = "0.0088"
$q3.Cells.Item(27, 8).Value = 5

$q3.Cells.Item(28, 1).Value = 26
$q3.Cells.Item(28, 2).Value = "002409"
$q3.Cells.Item(28, 3).Value = "华夏新活力灵活配置混合A"
$q3.Cells.Item(28, 4).Value = "0.15"
$q3.Cells.Item(28, 5).Value = "69.89"
$q3.Cells.Item(28, 6).Value = "3.96"
$q3.Cells.Item(28, 7).Value = "0.0059"
$q3.Cells.Item(28, 8).Value = 5

$q3.Cells.Item(29, 1).Value = 27
$q3.Cells.Item(29, 2).Value = "008842"
$q3.Cells.Item(29, 3).Value = "同泰远见灵活配置混合A"
$q3.Cells.Item(29, 4).Value = "0.18"
$q3.Cells.Item(29, 5).Value = "93.90"
$q3.Cells.Item(29, 6).Value = "2.90"
$q3.Cells.Item(29, 7).Value = "0.0052"
$q3.Cells.Item(29, 8).Value = 6

$q3.Cells.Item(30, 1).Value = 28
$q3.Cells.Item(30, 2).Value = "970046"
$q3.Cells.Item(30, 3).Value = "东海证券海睿健行灵活配置混合A"
$q3.Cells.Item(30, 4).Value = "0.16"
$q3.Cells.Item(30, 5).Value = "82.61"
$q3.Cells.Item(30, 6).Value = "3.01"
$q3.Cells.Item(30, 7).Value = "0.0048"
$q3.Cells.Item(30, 8).Value = 9

$q3.Cells.Item(31, 1).Value = 29
$q3.Cells.Item(31, 2).Value = "008890"
$q3.Cells.Item(31, 3).Value = "中邮价值优选一年定期开放灵活配置混合"
$q3.Cells.Item(31, 4).Value = "0.12"
$q3.Cells.Item(31, 5).Value = "61.34"
$q3.Cells.Item(31, 6).Value = "3.31"
$q3.Cells.Item(31, 7).Value = "0.0040"
$q3.Cells.Item(31, 8).Value = 7

$q3.Cells.Item(32, 1).Value = 30
$q3.Cells.Item(32, 2).Value = "014361"
$q3.Cells.Item(32, 3).Value = "红塔红土稳健添利混合C"
$q3.Cells.Item(32, 4).Value = "0.27"
$q3.Cells.Item(32, 5).Value = "37.80"
$q3.Cells.Item(32, 6).Value = "1.45"
$q3.Cells.Item(32, 7).Value = "0.0039"
$q3.Cells.Item(32, 8).Value = 7

$q3.Cells.Item(33, 1).Value = 31
$q3.Cells.Item(33, 2).Value = "015772"
$q3.Cells.Item(33, 3).Value = "银华互联网主题灵活配置混合C"
$q3.Cells.Item(33, 4).Value = "0.09"
$q3.Cells.Item(33, 5).Value = "89.52"
$q3.Cells.Item(33, 6).Value = "3.61"
$q3.Cells.Item(33, 7).Value = "0.0032"
$q3.Cells.Item(33, 8).Value = 10

$q3.Cells.Item(34, 1).Value = 32
$q3.Cells.Item(34, 2).Value = "970047"
$q3.Cells.Item(34, 3).Value = "东海证券海睿健行灵活配置混合B"
$q3.Cells.Item(34, 4).Value = "0.10"
$q3.Cells.Item(34, 5).Value = "82.61"
$q3.Cells.Item(34, 6).Value = "3.01"
$q3.Cells.Item(34, 7).Value = "0.0030"
$q3.Cells.Item(34, 8).Value = 9

$q3.Cells.Item(35, 1).Value = 33
$q3.Cells.Item(35, 2).Value = "015694"
$q3.Cells.Item(35, 3).Value = "瑞达策略优选混合A"
$q3.Cells.Item(35, 4).Value = "0.09"
$q3.Cells.Item(35, 5).Value = "67.87"
$q3.Cells.Item(35, 6).Value = "2.82"
$q3.Cells.Item(35, 7).Value = "0.0025"
$q3.Cells.Item(35, 8).Value = 4

$q3.Cells.Item(36, 1).Value = 34
$q3.Cells.Item(36, 2).Value = "014871"
$q3.Cells.Item(36, 3).Value = "摩根士丹利华鑫科技领先灵活配置混合C"
$q3.Cells.Item(36, 4).Value = "0.08"
$q3.Cells.Item(36, 5).Value = "94.13"
$q3.Cells.Item(36, 6).Value = "3.18"
$q3.Cells.Item(36, 7).Value = "0.0025"
$q3.Cells.Item(36, 8).Value = 9

$q3.Cells.Item(37, 1).Value = 35
$q3.Cells.Item(37, 2).Value = "005438"
$q3.Cells.Item(37, 3).Value = "易方达易百智能量化策略灵活配置混合C"
$q3.Cells.Item(37, 4).Value = "0.21"
$q3.Cells.Item(37, 5).Value = "94.77"
$q3.Cells.Item(37, 6).Value = "1.09"
$q3.Cells.Item(37, 7).Value = "0.0023"
$q3.Cells.Item(37, 8).Value = 5

$q3.Cells.Item(38, 1).Value = 36
$q3.Cells.Item(38, 2).Value = "008843"
$q3.Cells.Item(38, 3).Value = "同泰远见灵活配置混合C"
$q3.Cells.Item(38, 4).Value = "0.07"
$q3.Cells.Item(38, 5).Value = "93.90"
$q3.Cells.Item(38, 6).Value = "2.90"
$q3.Cells.Item(38, 7).Value = "0.0020"
$q3.Cells.Item(38, 8).Value = 6

$q3.Cells.Item(39, 1).Value = 37
$q3.Cells.Item(39, 2).Value = "004727"
$q3.Cells.Item(39, 3).Value = "先锋聚优灵活配置混合C"
$q3.Cells.Item(39, 4).Value = "0.03"
$q3.Cells.Item(39, 5).Value = "91.74"
$q3.Cells.Item(39, 6).Value = "2.58"
$q3.Cells.Item(39, 7).Value = "0.0008"
$q3.Cells.Item(39, 8).Value = 5

$q3.Cells.Item(40, 1).Value = 38
$q3.Cells.Item(40, 2).Value = "004726"
$q3.Cells.Item(40, 3).Value = "先锋聚优灵活配置混合A"
$q3.Cells.Item(40, 4).Value = "0.01"
$q3.Cells.Item(40, 5).Value = "91.74"
$q3.Cells.Item(40, 6).Value = "2.58"
$q3.Cells.Item(40, 7).Value = "0.0003"
$q3.Cells.Item(40, 8).Value = 5

$q3.Cells.Item(41, 1).Value = 39
$q3.Cells.Item(41, 2).Value = "166107"
$q3.Cells.Item(41, 3).Value = "信澳量化多因子混合（LOF）A"
$q3.Cells.Item(41, 4).Value = "0.05"
$q3.Cells.Item(41, 5).Value = "28.39"
$q3.Cells.Item(41, 6).Value = "0.40"
$q3.Cells.Item(41, 7).Value = "0.0002"
$q3.Cells.Item(41, 8).Value = 6

$q3.Cells.Item(42, 1).Value = 40
$q3.Cells.Item(42, 2).Value = "166108"
$q3.Cells.Item(42, 3).Value = "信澳量化多因子混合（LOF）C"
$q3.Cells.Item(42, 4).Value = "0.06"
$q3.Cells.Item(42, 5).Value = "28.39"
$q3.Cells.Item(42, 6).Value = "0.40"
$q3.Cells.Item(42, 7).Value = "0.0002"
$q3.Cells.Item(42, 8).Value = 6

$q3.Cells.Item(43, 1).Value = 41
$q3.Cells.Item(43, 2).Value = "015695"
$q3.Cells.Item(43, 3).Value = "瑞达策略优选混合C"
$q3.Cells.Item(43, 4).Value = "0.00"
$q3.Cells.Item(43, 5).Value = "67.87"
$q3.Cells.Item(43, 6).Value = "2.82"
$q3.Cells.Item(43, 7).NumberFormat = "General"
$q3.Cells.Item(43, 7).Value = 0
$q3.Cells.Item(43, 8).Value = 4

$q3.Cells.Item(44, 1).Value = 42
$q3.Cells.Item(44, 2).Value = "015147"
$q3.Cells.Item(44, 3).Value = "华夏兴融灵活配置混合（LOF）C"
$q3.Cells.Item(44, 4).Value = "0.00"
$q3.Cells.Item(44, 5).Value = "89.61"
$q3.Cells.Item(44, 6).Value = "2.24"
$q3.Cells.Item(44, 7).NumberFormat = "General"
$q3.Cells.Item(44, 7).Value = 0
$q3.Cells.Item(44, 8).Value = 10

$q3.Cells.Item(45, 1).Value = 43
$q3.Cells.Item(45, 2).Value = "002410"
$q3.Cells.Item(45, 3).Value = "华夏新活力灵活配置混合C"
$q3.Cells.Item(45, 4).Value = "0.00"
$q3.Cells.Item(45, 5).Value = "69.89"
$q3.Cells.Item(45, 6).Value = "3.96"
$q3.Cells.Item(45, 7).NumberFormat = "General"
$q3.Cells.Item(45, 7).Value = 0
$q3.Cells.Item(45, 8).Value = 5

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row for 2022-Q3 above the
#    existing quarters, shifting the older rows down by one and renumbering
#    the running ID column (A).
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

for ($r = 6; $r -ge 2; $r--) {
    $dest = $r + 1
    $summary.Cells.Item($dest, 1).Value = $r - 1
    $summary.Cells.Item($dest, 2).Value = $summary.Cells.Item($r, 2).Value2
    $summary.Cells.Item($dest, 3).Value = $summary.Cells.Item($r, 3).Value2
    $summary.Cells.Item($dest, 4).Value = $summary.Cells.Item($r, 4).Value2
}

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 44
$summary.Cells.Item(2, 4).Value = 3.56

# Row 7 is brand new (shifted down from nothing before it) - give column A
# the same bold/centered/bordered style used by the other ID cells.
$a7 = $summary.Cells.Item(7, 1)
$a7.Font.Bold = $true
$a7.HorizontalAlignment = -4108
$a7.VerticalAlignment = -4160
$a7.Borders.LineStyle = 1
